$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the existing "Status" header (column M) to "Order Status".
$ws.Range("M1").Value = "Order Status"

# 2. Append a new "Paid Status" column right after it (column N), matching
#    the header formatting used by the other plain header cells.
$ws.Range("A1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N1").Value = "Paid Status"

# 3. Insert a new column before "Amount" (L) for "Return Items". This shifts
#    Amount L->M, Order Status M->N and Paid Status N->O.
$ws.Columns("L:L").Insert()
$ws.Range("L1").Value = "Return Items"

# 4. Insert four blank rows above the header row, pushing it from row 1 to row 5.
$ws.Rows("1:4").Insert()

# 5. Add the "Sales Agent" / "Date" labels in the newly freed rows 1-2.
$ws.Range("A1").Value = "Sales Agent"
$ws.Range("A1").Font.Bold = $true

$ws.Range("A2").Value = "Date"
$ws.Range("A2").Font.Bold = $true

# 6. Match the active selection recorded in the saved workbook.
$ws.Range("B8").Select()
